# "Generate Report for Archive"
# The localization status for the two tracked files flips from
# "Ready for handoff" to "In Translation" everywhere it is shown
# (the Overview sheet's zh-cn/de-de columns, and the Status column
# on each per-locale detail sheet). Excel's column AutoFit then
# narrows the now-shorter "Status" columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn (col E) / de-de (col F) status cells for both rows.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn detail sheet: Status column (C) for both rows.
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

# de-de detail sheet: Status column (C) for both rows.
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# The Status columns shrink now that "In Translation" is shorter than
# "Ready for handoff" - resize them to match the new content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
